$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Content edits -----------------------------------------------------

# Week 8 ("Shiny 1 & 2") description gets a new reminder prefixed on to it,
# and its slide deck ("w8") is now posted.
$ws.Range("E9").Value = "Please clone [the course repo of example apps](https://github.com/datalorax/shiny-app-examples) before class.We will start by introducing the very basics of shiny - the user interface (UI) and the server. We will work together to create a basic shiny application, modifying the default template to use ggplot2. Shiny dashboards and different layout options will also be discussed."
$ws.Range("D9").Value = "w8"

# Week 9 ("Shiny 3 & review") slide deck ("w9") is now posted.
$ws.Range("D10").Value = "w9"

# The week-8 row grew taller once the extra reminder sentence was added to
# its (wrapped) description cell.
$ws.Rows.Item(9).RowHeight = 119

# --- View / selection state ---------------------------------------------
# The author scrolled the sheet down to row 9 and over to column D (just
# past the frozen A:C columns) and left the selection on F10.
$win = $excel.ActiveWindow
$ws.Range("F10").Select()

$leftPane = $win.Panes.Item(1)
$rightPane = $win.Panes.Item(2)
$leftPane.ScrollRow = 9
$leftPane.ScrollColumn = 4
$rightPane.ScrollRow = 9
$rightPane.ScrollColumn = 4
